$d = $word.ActiveDocument

# Use Find.Execute to *locate* text (no replacement text passed), then set
# Range.Text directly. Doing the substitution this way avoids the
# "smart quotes" autocorrection that Find/Replace's own Replace-With
# argument applies to straight apostrophes/quotes, which would otherwise
# diverge from the target text that uses plain ASCII punctuation.
function Replace-Text($find, $replace) {
    $count = 0
    $r = $d.Content
    while ($r.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, "", 0)) {
        $r.Text = $replace
        $count++
        if ($count -gt 20) { break }
        $r = $d.Range($r.End, $d.Content.End)
    }
}

# [Music] -> [Musique]  (occurs twice; ReplaceAll handles both)
Replace-Text "[Music]" "[Musique]"

Replace-Text "après longtemps. After some" "après longtemps. Après quelques"

Replace-Text "chatting, Phil says he has three children, then" "discussions, Phil dit qu'il a trois enfants, puis"

Replace-Text "Mike, astonished, asks: 'How old are they?' Fil," "Mike, étonné, demande : « Quel âge ont-ils ? » Fil,"

Replace-Text "being a playful mathematician, answers" "qui est un mathématicien ludique, lui répond"

Replace-Text "'You tell me! I'll give you a hint: if you" "« Dis-le moi ! Je vais te donner un indice :"

Replace-Text "multiply the three ages together you" "les trois âges multipliés donnent"

Replace-Text "get 36.' Mike takes sometimes to think" "36. » Mike prend un moment à réfléchir"

Replace-Text "and says: 'I'm sorry Fil, but I do need" "et dit: « Je suis désolé Fil, mais j'ai besoin"

Replace-Text "another hint. So Fil tells Mike:" "d'un autre indice. » Alors, Fil dit à Mike :"

Replace-Text "'Yes, sure, here it is: if you had up to" "« Oui, bien sûr, voilà: si tu additionnes les"

Replace-Text "three ages you get the number of math" "trois âges, tu obtiens le nombre d'articles"

Replace-Text "papers we publish together. Do you remember it?'" "mathématiques que nous avons publiés ensemble. Tu te souviens ? »"

Replace-Text "'Yes I do remember How many, but still" "« Oui je m'en souviens mais"

Replace-Text "I do not have enough information! I need" "je n'ai toujours pas assez d'informations ! J'ai besoin d'"

Replace-Text "at least one more.' Fil says: 'Yes don't" "au moins un de plus. » Fil dit : « Oui, pas"

Replace-Text "worry but this is the last one:" "de soucis, mais c'est le dernier :"

Replace-Text "The youngest one has blues eyes.' And" "« Le plus jeune a les yeux bleus. » Et"

Replace-Text "suddenly Mike gets the answer. You" "soudainement Mike obtient la réponse. Toi, tu"

Replace-Text "hear the conversation but you don't know" "entends la conversation mais tu ne sais pas"

Replace-Text "how many papers they published together." "combien d'articles ils ont publié ensemble."

Replace-Text "However, you do want to know the ages of" "Toutefois, tu souhaites connaître l'âge des"

Replace-Text "the three children. Can you figure them" "trois enfants. Peux-tu les trouver ?"

Replace-Text "out?" " "
